$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.363.87"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.380.95"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.14"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.69"
$ws.Range("E6").Value = "  +9.53%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.381.84"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.62"
$ws.Range("E10").Value = "  +5.19%  "
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  +5.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.947.71"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +3.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.372.25"
$ws.Range("E16").Value = "  +1.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.28"
$ws.Range("E17").Value = "  +3.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.359.42"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.00"
$ws.Range("E19").Value = "  +6.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.48"
$ws.Range("E20").Value = "  +4.59%  "
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.16"
$ws.Range("E22").Value = "  +9.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.576"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.508.06"
$ws.Range("E24").Value = "  +1.59%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.86"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  +11.48%  "
$ws.Range("E28").Value = "  +15.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.84"
$ws.Range("E29").Value = "  +9.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.156"
$ws.Range("E32").Value = "  +6.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.14"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.407.62"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.51"
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.58"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.98"
$ws.Range("E38").Value = "  +4.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +5.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.17"
$ws.Range("E40").Value = "  +0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0805"
$ws.Range("E41").Value = "  +7.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.43"
$ws.Range("E43").Value = "  +5.15%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.21"
$ws.Range("E44").Value = "  +9.17%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.55"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.763"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.66"
$ws.Range("E47").Value = "  +7.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.46"
$ws.Range("E48").Value = "  +3.77%  "
$ws.Range("E49").Value = "  +5.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.24"
$ws.Range("E50").Value = "  +12.96%  "
$ws.Range("E51").Value = "  +13.29%  "
